$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three new columns (D, E, F) before the existing "Terms Typically Offered" column,
# which shifts it from D to G.
$ws.Columns("D:F").Insert()

# New header row for the inserted columns
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# Per-row data: row number, new Prerequisites (C) text, Corequisites (D),
# Concurrent (E), Recommended (F), Terms Typically Offered (G)
$rows = @(
  @{ Row=2; C="NA"; D="NA"; E="NA"; F="NA"; G="TBD" },
  @{ Row=3; C="NA"; D="NA"; E="NA"; F="NA"; G="F, SP" },
  @{ Row=4; C="Consent of instructor."; D="NA"; E="NA"; F="NA"; G="TBD" },
  @{ Row=5; C="DSCI 230 and DSCI 231."; D="NA"; E="NA"; F="NA"; G="F" },
  @{ Row=6; C="NA"; D="NA"; E="NA"; F="NA"; G="F, W, SP" },
  @{ Row=7; C="NA"; D="NA"; E="NA"; F="NA"; G="F, SP" },
  @{ Row=8; C="DSCI 231, or FSN 125, or FSN 230."; D="NA"; E="NA"; F="NA"; G="W" },
  @{ Row=9; C="DSCI 230."; D="NA"; E="NA"; F="NA"; G="SP" },
  @{ Row=10; C="Open to undergraduate students and consent of instructor."; D="NA"; E="NA"; F="NA"; G="TBD" },
  @{ Row=11; C="ASCI 220 with a grade of C- or better or consent of instructor."; D="NA"; E="NA"; F="NA"; G="F, SP" },
  @{ Row=12; C="ASCI 220; DSCI 230; and BIO 111 or BIO 161."; D="NA"; E="NA"; F="NA"; G="F" },
  @{ Row=13; C="ASCI 229 or DSCI 230."; D="NA"; E="NA"; F="NA"; G="F" },
  @{ Row=14; C="ASCI 220 and DSCI 230."; D="NA"; E="NA"; F="NA"; G="W" },
  @{ Row=15; C="Consent of internship instructor."; D="NA"; E="NA"; F="NA"; G="TBD" },
  @{ Row=16; C="Consent of instructor."; D="NA"; E="NA"; F="NA"; G="TBD" },
  @{ Row=17; C="Junior standing."; D="NA"; E="NA"; F="CHEM 312."; G="SP " },
  @{ Row=18; C="DSCI 444 or MCRO 421."; D="NA"; E="NA"; F="NA"; G="F" },
  @{ Row=19; C="ASCI 229; and ASCI 355 or DSCI 301."; D="NA"; E="NA"; F="NA"; G="F" },
  @{ Row=20; C="DSCI 333."; D="NA"; E="NA"; F="NA"; G="W" },
  @{ Row=21; C="DSCI 241, BIO 111 or higher, STAT 130 or higher."; D="NA"; E="NA"; F="NA"; G="SP" },
  @{ Row=22; C="DSCI 333."; D="NA"; E="NA"; F="NA"; G="W" },
  @{ Row=23; C="DSCI 233 or FSN 204."; D="NA"; E="NA"; F="NA"; G="F" },
  @{ Row=24; C="DSCI 231; MCRO 221 or MCRO 224; and STAT 218."; D="NA"; E="NA"; F="CHEM 313."; G="W " },
  @{ Row=25; C="DSCI 233 or FSN 204."; D="NA"; E="NA"; F="NA"; G="SP" },
  @{ Row=26; C="DSCI 233; and MCRO 221 or MCRO 224; and STAT 130 or STAT 218; or graduate standing."; D="NA"; E="NA"; F="NA"; G="TBD" },
  @{ Row=27; C="Junior standing."; D="NA"; E="NA"; F="NA"; G="TBD" },
  @{ Row=28; C="Consent of instructor."; D="NA"; E="NA"; F="NA"; G="TBD" },
  @{ Row=29; C="Consent of instructor."; D="NA"; E="NA"; F="NA"; G="TBD" },
  @{ Row=30; C="Consent of department head, graduate advisor and supervising faculty member."; D="NA"; E="NA"; F="NA"; G="TBD" },
  @{ Row=31; C="CHEM 212 or CHEM 312; admission to Master of Professional Studies program in Dairy Products Technology."; D="NA"; E="NA"; F="Differential and Integral Calculus."; G="F " },
  @{ Row=32; C="DSCI 501; admission to Master of Professional Studies program in Dairy Products Technology."; D="NA"; E="NA"; F="NA"; G="W" },
  @{ Row=33; C="Admission to Master of Professional Studies program in Dairy Products Technology."; D="NA"; E="NA"; F="NA"; G="F" },
  @{ Row=34; C="Admission to Master of Professional Studies program in Dairy Products Technology."; D="NA"; E="DSCI 540."; F="NA"; G="F " },
  @{ Row=35; C="DSCI 520; admission to Master of Professional Studies program in Dairy Products Technology."; D="NA"; E="NA"; F="NA"; G="W" },
  @{ Row=36; C="DSCI 520 and admission to Master of Professional Studies program in Dairy Products Technology."; D="NA"; E="NA"; F="NA"; G="SP" },
  @{ Row=37; C="DSCI 501 and admission to Master of Professional Studies program in Dairy Products Technology."; D="NA"; E="NA"; F="NA"; G="SP" },
  @{ Row=38; C="Admission to Master of Professional Studies program in Dairy Products Technology and consent of instructor."; D="NA"; E="NA"; F="NA"; G="SU" },
  @{ Row=39; C="Admission to Master of Professional Studies program in Dairy Products Technology."; D="NA"; E="DSCI 520."; F="MCRO 221 or equivalent."; G="F  " },
  @{ Row=40; C="Admission to Master of Professional Studies program in Dairy Products Technology."; D="NA"; E="NA"; F="NA"; G="W" },
  @{ Row=41; C="Senior or graduate standing and approval of instructor."; D="NA"; E="NA"; F="NA"; G="TBD" },
  @{ Row=42; C="DSCI 520; admission to Master of Professional Studies program in Dairy Science Technology."; D="NA"; E="NA"; F="NA"; G="SP" },
  @{ Row=43; C="Graduate standing or consent of instructor."; D="NA"; E="NA"; F="NA"; G="TBD" },
  @{ Row=44; C="Consent of instructor."; D="NA"; E="NA"; F="NA"; G="TBD" },
  @{ Row=45; C="Admission to Master of Professional Studies program in Dairy Products Technology."; D="NA"; E="NA"; F="NA"; G="F, W, SP" },
  @{ Row=46; C="Admission to Master of Professional Studies program in Dairy Products Technology."; D="NA"; E="NA"; F="NA"; G="W" },
  @{ Row=47; C="Graduate standing and consent of instructor."; D="NA"; E="NA"; F="NA"; G="TBD" },
  @{ Row=48; C="Graduate standing."; D="NA"; E="NA"; F="NA"; G="TBD" }
)

foreach ($item in $rows) {
  $r = $item.Row
  $ws.Range("C$r").Value = $item.C
  $ws.Range("D$r").Value = $item.D
  $ws.Range("E$r").Value = $item.E
  $ws.Range("F$r").Value = $item.F
  $ws.Range("G$r").Value = $item.G
}

